# Appends a new 03/06/2024 log entry to the end of the development log,
# mirroring the structure already used for every other dated entry:
# a blank paragraph, then the date, then one or more body paragraphs,
# each separated by a blank paragraph.

$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

function New-ParaXml {
    param($Runs)
    if (-not $Runs) {
        return '<w:p xmlns:w="' + $wNs + '"/>'
    }
    $inner = ''
    foreach ($run in $Runs) {
        $pb = ''
        if ($run.pageBreak) {
            $pb = '<w:lastRenderedPageBreak/>'
        }
        $attrs = ''
        if ($run.preserve) {
            $attrs = ' xml:space="preserve"'
        }
        $inner += '<w:r>' + $pb + '<w:t' + $attrs + '>' + $run.text + '</w:t></w:r>'
    }
    return '<w:p xmlns:w="' + $wNs + '">' + $inner + '</w:p>'
}

$paraXmls = @(
    # blank separator line before the new entry
    (New-ParaXml $null),

    # date heading
    (New-ParaXml @(
        @{ text = '03/06/2024' }
    )),

    # manufacturing capabilities paragraph
    (New-ParaXml @(
        @{ text = 'Since my last update I improved my manufacturing capabilities through '; preserve = $true },
        @{ text = 'lobbying and work at school. We now have an operational lab furnace which can reach temperatures of 1850 degrees F. This opens the possibility of casting components of the frame in the machine shop. This combined with improved access makes making various components much more feasible. With this, I plan to rough cast aluminum parts for the pistons, frame base assemblies, and some components of the drive train and then machine them to proper dimensions using the schools CNC router/lathe as needed. To connect the bases I plan to cut the stainless steel I acquired while working at Ironheart welding using a plasma cutter (also at ironheart welding). Other components will be 3D printed such as valve actuators, UI and electronics mounts, pipe supports, and other non-load-bearing components of the design.' }
    )),

    (New-ParaXml $null),

    # machine shop access paragraph
    (New-ParaXml @(
        @{ text = 'As mentioned above I now have better access to the school machine shop' },
        @{ text = ', meaning I can now reasonably get access for an entire day and on weekends when I do not have obligations such as class' },
        @{ text = '. This gives me room to make' },
        @{ text = ' a clear and actionable plan for building the washing machine. Starting this weekend, I will be making molds for the rough casts made for the bases'; preserve = $true },
        @{ text = ', which will just be an open, rectangular mold welded from scrap steel in the machine shop.' }
    )),

    (New-ParaXml $null),

    # electronics / professor collaboration paragraph (note the lastRenderedPageBreak on the first run)
    (New-ParaXml @(
        @{ text = 'With this, I have also began working with one of m'; pageBreak = $true },
        @{ text = 'y former professors to start taking actionable steps towards the electronics and other '; preserve = $true },
        @{ text = 'similar portions of the project' },
        @{ text = '.' }
    )),

    (New-ParaXml $null),

    # closing paragraph
    (New-ParaXml @(
        @{ text = 'For the electronics, I still need to approximate the load for a required motor size, but once this has been done, all the needed electronic components for the project will be ready for building.' }
    ))
)

# Walk a collapsed range forward, inserting a fresh paragraph each time and
# then overwriting that paragraph's contents with the exact OOXML we want
# (this keeps truly-blank separators as empty <w:p/> elements rather than
# paragraphs holding a stray empty run).
$rng = $d.Paragraphs.Last.Range
$rng.Collapse(0)

foreach ($paraXml in $paraXmls) {
    $rng.InsertParagraphAfter()
    $allParas = $d.Paragraphs
    $newPara = $allParas.Item($allParas.Count)
    $newPara.Range.InsertXML($paraXml) | Out-Null
    $rng = $newPara.Range
    $rng.Collapse(0)
}

Write-Output "Inserted $($paraXmls.Count) paragraphs; document now has $($d.Paragraphs.Count) paragraphs."
